$wb = $excel.ActiveWorkbook

# The "Croatia" sheet is a copy of the existing "Turkey" sheet, inserted
# right after it, with updated market name / reference values.
$turkey = $wb.Worksheets.Item("Turkey")
$turkey.Copy($null, $turkey)

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3191/T2423"

# Restore Turkey's own selection to span the whole sheet (as it ends up
# after Excel re-saves it once it is no longer the active tab), then make
# Croatia the active sheet/tab with B4 selected.
$turkey.Range("A1:XFD1048576").Select() | Out-Null
$croatia.Activate() | Out-Null
$croatia.Range("B4").Select() | Out-Null
